$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new text looks like a plain number (e.g. "300.23") would be
# auto-converted to a numeric value by Excel when set via .Value. For those cells we
# temporarily force a Text number format, assign the literal string, then restore the
# cell style to Normal so the final style matches the original (unstyled) cell.

$ws.Range("D2").Value = '45.730.86'
$ws.Range("E2").Value = '  -2.61%  '

$ws.Range("D3").Value = '2.352.96'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0800'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.02%  '

$ws.Range("E13").Value = '  -1.40%  '

$ws.Range("D14").Value = '2.710.97'
$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("D15").Value = '2.361.28'
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.813'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.93%  '

$ws.Range("D18").Value = '45.693.19'
$ws.Range("E18").Value = '  -2.42%  '

$ws.Range("D19").Value = '0.0₃0975'
$ws.Range("E19").Value = '  +2.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '40.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +17.51%  '

$ws.Range("E32").Value = '  +4.77%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.35%  '

$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '145.28'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0774'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.116'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0299'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("D44").Value = '1.858.24'
$ws.Range("E44").Value = '  +2.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.185'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '70.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.99%  '

$ws.Range("D49").Value = '2.582.51'
$ws.Range("E49").Value = '  +0.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '96.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.75%  '

